# Terminal Hortofrutícola Agro Chillán - Repollo
# Insert a new weekly record at row 52, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52 (pushes existing rows 52..162 down to 53..163)
$ws.Rows.Item(52).Insert()

# Populate the new row 52 with the new weekly observation
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 44544
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 100112006
$ws.Range("G52").Value = "Repollo"
$ws.Range("H52").Value = "Crespo record"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 1000
$ws.Range("K52").Value = 500
$ws.Range("L52").Value = 600
$ws.Range("M52").Value = 550
$ws.Range("N52").Value = "$/unidad"
$ws.Range("O52").Value = "Provincia de Diguillín"
$ws.Range("P52").Value = 550
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = "Hortaliza"
